$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Sheet1 edits
# ---------------------------------------------------------------------------

# Drop the old hyperlinks (email / password mailto's + the cart-page link);
# a single new hyperlink is (re)created on B5 below.
$ws.Hyperlinks.Delete()

# Row 5 used to be "waitfortext / Continue"; it becomes a "goto" of the
# Blue Nile product page, with wait timers, and a hyperlink + Hyperlink
# style on the URL cell (matching B10/B15's existing look).
$blueNileUrl = "https://www.bluenile.com/jewelry/necklaces/lab-grown-diamond-cushion-cut-solitaire-pendant-in-14k-white-gold-1-2-ct-tw-f-g-vs2-si1-item-202314"
$ws.Range("A5").Value = "goto"
$ws.Range("B5").Value = $blueNileUrl
$ws.Range("D5").Value = 1000
$ws.Range("E5").Value = 9000
$ws.Hyperlinks.Add($ws.Range("B5"), $blueNileUrl) | Out-Null
$ws.Range("B10").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 6 used to be the "type" of the e-mail address (with a hyperlinked
# value cell); it becomes a "scroll" step waiting for "Ships by".
$ws.Range("A6").Value = "scroll"
$ws.Range("B6").Value = "Ships by"
$ws.Range("C6").Clear()
$ws.Range("D6").Value = 1000
$ws.Range("E6").Value = 5000

# Row 7 used to be "click / Continue button on page"; it becomes an
# "ai_click" on the "ADD TO CART button".
$ws.Range("A7").Value = "ai_click"
$ws.Range("B7").Value = "ADD TO CART button"
$ws.Range("D7").Value = 1000
$ws.Range("E7").Value = 2000

# Row 8 used to be "type / password input field in span / Welcome@123456 /
# 1000 / 1000"; it collapses down to a single "clickto" action.
$ws.Range("A8").Value = "clickto"
$ws.Range("B8:E8").Clear()

# Row 9 ("click / signin button on page / 3000") is removed entirely.
$ws.Rows("9:9").Clear()

# Row 10 used to be "goto / cart url / 5000"; only the (empty, still
# hyperlink-styled) B10 cell survives.
$ws.Range("A10").Clear()
$ws.Range("B10").ClearContents()
$ws.Range("E10").Clear()

# Rows 11-13 (checkvisible / click / assert steps) are removed entirely.
$ws.Rows("11:13").Clear()

# The active cell moves to A8.
$ws.Range("A8").Select() | Out-Null

# ---------------------------------------------------------------------------
# New "Sheet2" (TID / Execution results) is inserted right after Sheet1.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "TID"
$ws2.Range("B1").Value = "Execution"
$ws2.Range("A2").Value = "T1"
$ws2.Range("B2").Value = "Y"
$ws2.Range("A3").Value = "T2"
$ws2.Range("B3").Value = "N"

# Sheet1 stays the active/selected tab.
$ws.Activate()
